{"js": "// Replace the 25 division-expression runs in the table (in document order)\n// with their new values, leaving the blank \"answer\" rows untouched.\nconst newValues = [\n  \"411\u00f75=\", \"910\u00f73=\", \"328\u00f79=\", \"476\u00f78=\", \"898\u00f73=\",\n  \"727\u00f76=\", \"587\u00f72=\", \"290\u00f78=\", \"607\u00f73=\", \"136\u00f79=\",\n  \"218\u00f72=\", \"222\u00f79=\", \"397\u00f76=\", \"340\u00f77=\", \"808\u00f73=\",\n  \"277\u00f75=\", \"194\u00f77=\", \"266\u00f75=\", \"868\u00f75=\", \"323\u00f75=\",\n  \"998\u00f78=\", \"282\u00f79=\", \"719\u00f75=\", \"335\u00f75=\", \"164\u00f73=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet valueIndex = 0;\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (const cell of cells.items) {\n      cell.body.load(\"text\");\n      await context.sync();\n\n      const text = cell.body.text.trim();\n      if (text.length === 0) {\n        continue; // blank answer cell, leave as-is\n      }\n      if (valueIndex >= newValues.length) {\n        continue;\n      }\n\n      const paragraphs = cell.body.paragraphs;\n      paragraphs.load(\"items\");\n      await context.sync();\n\n      const paragraph = paragraphs.items[0];\n      paragraph.insertText(newValues[valueIndex], Word.InsertLocation.replace);\n      valueIndex++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-expression cells in the table (in document order)\n# with their new values, leaving the blank \"answer\" rows untouched.\n$newValues = @(\n    \"411\u00f75=\", \"910\u00f73=\", \"328\u00f79=\", \"476\u00f78=\", \"898\u00f73=\",\n    \"727\u00f76=\", \"587\u00f72=\", \"290\u00f78=\", \"607\u00f73=\", \"136\u00f79=\",\n    \"218\u00f72=\", \"222\u00f79=\", \"397\u00f76=\", \"340\u00f77=\", \"808\u00f73=\",\n    \"277\u00f75=\", \"194\u00f77=\", \"266\u00f75=\", \"868\u00f75=\", \"323\u00f75=\",\n    \"998\u00f78=\", \"282\u00f79=\", \"719\u00f75=\", \"335\u00f75=\", \"164\u00f73=\"\n)\n\n$d = $word.ActiveDocument\n$valueIndex = 0\n\nforeach ($t in $d.Tables) {\n    foreach ($cell in $t.Range.Cells) {\n        # Cell.Range.Text always ends with the paragraph mark (\\r) plus the\n        # cell mark (\\x07); a truly empty cell's text is just those two\n        # characters, so strip them off before checking for content.\n        $raw = $cell.Range.Text\n        $content = $raw.TrimEnd([char]7).TrimEnd([char]13)\n        if ($content.Length -eq 0) {\n            continue\n        }\n        if ($valueIndex -ge $newValues.Count) {\n            continue\n        }\n        $cell.Range.Text = $newValues[$valueIndex]\n        $valueIndex++\n    }\n}\n"}
